$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value for the Price (D) / Volume(1h) (E) columns refreshed by this run.
$updates = [ordered]@{
    'D2' = '63.502.58'
    'E2' = '  +0.05%  '
    'D3' = '2.639.73'
    'E3' = '  -0.34%  '
    'E4' = '  -0.01%  '
    'D5' = '604.32'
    'E5' = '  +2.04%  '
    'D6' = '146.47'
    'E6' = '  +1.98%  '
    'E7' = '  +0.00%  '
    'D8' = '0.588'
    'E8' = '  +0.33%  '
    'E9' = '  +2.23%  '
    'D10' = '5.59'
    'E10' = '  -0.37%  '
    'D11' = '0.369'
    'E11' = '  +4.32%  '
    'E12' = '  -0.35%  '
    'D13' = '27.69'
    'E13' = '  +1.08%  '
    'D14' = '3.112.71'
    'E14' = '  -0.41%  '
    'D15' = '63.347.82'
    'E15' = '  -0.07%  '
    'E16' = '  +1.50%  '
    'D17' = '2.638.22'
    'E17' = '  -0.42%  '
    'E18' = '  +0.73%  '
    'E19' = '  +4.76%  '
    'D20' = '344.09'
    'E20' = '  +1.14%  '
    'D21' = '6.94'
    'E21' = '  +2.90%  '
    'E22' = '  -0.19%  '
    'D23' = '5.56'
    'E23' = '  -3.74%  '
    'D24' = '66.62'
    'E24' = '  -0.89%  '
    'E25' = '  +1.97%  '
    'D26' = '587.79'
    'E26' = '  +8.98%  '
    'D27' = '9.08'
    'E27' = '  +7.85%  '
    'D28' = '1.54'
    'E28' = '  +0.56%  '
    'E29' = '  -1.79%  '
    'D30' = '7.98'
    'E30' = '  +2.12%  '
    'D31' = '1.00'
    'E31' = '  -0.19%  '
    'E32' = '  +3.52%  '
    'D33' = '1.75'
    'E33' = '  -2.99%  '
    'D34' = '0.0₃0831'
    'E34' = '  +3.10%  '
    'D35' = '5.22'
    'E35' = '  +7.44%  '
    'D36' = '166.78'
    'E36' = '  -4.59%  '
    'D37' = '0.406'
    'E37' = '  +1.06%  '
    'E38' = '  -0.13%  '
    'E39' = '  +8.20%  '
    'D40' = '19.11'
    'E40' = '  +0.21%  '
    'E41' = '  +0.04%  '
    'E42' = '  -2.45%  '
    'E43' = '  +1.61%  '
    'D44' = '22.19'
    'E44' = '  +0.04%  '
    'D45' = '0.0570'
    'E45' = '  +1.22%  '
    'D46' = '0.629'
    'E46' = '  -0.29%  '
    'E47' = '  +3.63%  '
    'D48' = '0.0963'
    'E48' = '  +0.24%  '
    'D49' = '1.92'
    'E49' = '  +13.14%  '
    'D50' = '18.79'
    'E50' = '  +0.74%  '
    'E51' = '  +1.93%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    # Price cells that look like plain numbers (e.g. "604.32") must stay text,
    # matching the sheet's original inlineStr storage (no trailing-zero loss).
    if ($cellRef.StartsWith("D")) {
        $looksNumeric = $value -match '^[+-]?\d+(\.\d+)?$'
        if ($looksNumeric) {
            $range.NumberFormat = "@"
        }
    }
    $range.Value = $value
}
